# Update Name of Algo
# Applies corrected numeric values to specific cells in Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B3").Value = 5.619000000000001
$ws.Range("B21").Value = 9.379000000000001
$ws.Range("B23").Value = 7.848999999999999
$ws.Range("B25").Value = 6.493
$ws.Range("E27").Value = 16.531
$ws.Range("E31").Value = 16.62
$ws.Range("E39").Value = 16.401
$ws.Range("E48").Value = 17.252
$ws.Range("E51").Value = 16.617
$ws.Range("E52").Value = 16.543
$ws.Range("B53").Value = 6.103999999999999
$ws.Range("E55").Value = 16.416
$ws.Range("E56").Value = 16.276
$ws.Range("B57").Value = 5.090999999999999
$ws.Range("E57").Value = 16.553
$ws.Range("B59").Value = 4.435
$ws.Range("B69").Value = 5.404
$ws.Range("E73").Value = 16.791
$ws.Range("B79").Value = 5.488
$ws.Range("B83").Value = 5.702000000000001
$ws.Range("E89").Value = 17.199
$ws.Range("E90").Value = 16.638
$ws.Range("B93").Value = 5.608
